# Applies the "Logs de prueba" update:
#  - D17 label: "Conectar con estacion" -> "Conectar con estacion*"
#  - F17:G17 becomes a merged "PDTE" cell, matching the style/format of F15:G15
#  - G11 and G13 are marked as Pass (green fill), matching the style of the
#    other "Pass" cells in column G
#  - Selection moves to J21

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the "Conectar con estacion" row label ---
$ws.Range("D17").Value = "Conectar con estacion*"

# --- Turn F17:G17 into a "PDTE" cell styled like F15:G15 ---
$ws.Range("F15:G15").Copy()
$ws.Range("F17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("F17:G17").Merge()
$ws.Range("F17").Value = "PDTE"

# --- Mark G11 and G13 as Pass (green), matching other Pass cells in column G ---
$ws.Range("G11").Interior.Color = $ws.Range("G9").Interior.Color
$ws.Range("G13").Interior.Color = $ws.Range("G9").Interior.Color

# --- Update the active selection ---
$ws.Range("J21").Select()
